$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new backlog item: "Padronizar CSS" (row 19) -----------------------
# Values first.
$ws.Range("D19").Value = "Padronizar CSS"
$ws.Range("E19").Value = "Padronização do CSS do site para facilitar a leitura do codigo"
$ws.Range("F19").Value = "Desejavel"
$ws.Range("G19").Value = "P"
$ws.Range("H19").Value = 5
$ws.Range("I19").Value = 2

# Mirror formatting (fill/border/alignment/number format) of the row above
# it (row 18), which uses the same alternating-stripe style.
$ws.Range("D18:I18").Copy()
$ws.Range("D19:I19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the row height used by the other wrapped-text rows in the table.
$ws.Rows.Item(19).RowHeight = 28.8

# --- Update the saved cursor/selection position -----------------------------
$ws.Range("L11").Select() | Out-Null
